$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd student email in C2: "pasne.d@husky.neu.edu" -> "panse.d@husky.neu.edu"
$ws.Range("C2").Value = "panse.d@husky.neu.edu"

# Leave the cursor/selection on the cell that was just edited
$ws.Range("C2").Select()
